# Add season-record columns (Wins / Losses / Ties) to the BAL_1994 sheet.
# Mirrors the author's fix: the original scrape only pulled team/player
# statistics, not the season W-L-T record, so three new columns are
# appended after the existing data (columns AD, AE, AF) with the header
# row labelled "Wins", "Losses", "Ties" and every data row (2-34) filled
# with the team's 1994 season record (63-49-0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - use the same bold/bordered header style as the
# neighbouring header cells (e.g. AC1) by copying it onto the new cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$wins = 63
$losses = 49
$ties = 0

for ($row = 2; $row -le 34; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
